$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.45642352104187
$ws.Range("B1").Value = 1.926920652389526
$ws.Range("C1").Value = 2.056244373321533
$ws.Range("D1").Value = 1.632256388664246
$ws.Range("E1").Value = 1.428162217140198
